$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translatable_Site_labels")
$ws.Activate()

# --- Row 41: "Select list type" -> "Import type" --------------------------
$ws.Range("C41").Value = "Import type"

# --- New rows 42-45: additional import-related labels ----------------------
$ws.Range("B42").Value = "import_preview_header"
$ws.Range("C42").Value = "Preview"
$ws.Range("D42").Value = "Import"
$ws.Range("E42").Formula = "=CONCAT(`"`",B42,`" : '`",C42,`"',`")"

$ws.Range("B43").Value = "import_upload_data"
$ws.Range("C43").Value = "Upload data"
$ws.Range("D43").Value = "Import"
$ws.Range("E43").Formula = "=CONCAT(`"`",B43,`" : '`",C43,`"',`")"

$ws.Range("B44").Value = "import_preview_label"
$ws.Range("C44").Value = "Please change column names using the dropdowns"
$ws.Range("D44").Value = "Import"
$ws.Range("E44").Formula = "=CONCAT(`"`",B44,`" : '`",C44,`"',`")"

$ws.Range("B45").Value = "import_push_data"
$ws.Range("C45").Value = "Push data to database"
$ws.Range("D45").Value = "Import"
$ws.Range("E45").Formula = "=CONCAT(`"`",B45,`" : '`",C45,`"',`")"

# --- View state: scroll + selection to match the saved workbook view -------
$win = $wb.Windows.Item(1)
$win.ScrollRow = 16
$win.ScrollColumn = 1
$ws.Range("B46").Select() | Out-Null
